# Generate Report for Handoff
#
# This re-runs the localization-status report generation for the rows that
# are currently "Ready for handoff" with a handoff priority of "ht".
# Their priority flips from the (stale) "ht" handoff-type marker to "mt"
# (the marker used once the handoff file has actually been generated/sent),
# and the handoff timestamp columns are refreshed to the moment the report
# was (re)generated. The Overview sheet's "Latest HO Xliff Generate Date"
# column mirrors the newest per-file handoff timestamp from the language
# sheets, so it is refreshed for the same rows.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Rows (2-based, header is row 1) on the language sheets that currently have
# Priority = "ht" and Status = "Ready for handoff".
$rows = @(7, 9, 10, 11, 12, 13, 14, 16)

foreach ($r in $rows) {
    $zhcn.Cells.Item($r, 5).Value = "mt"              # Priority
    $zhcn.Cells.Item($r, 8).Value = "2016-08-13 06:26:55"   # Latest Handoff Datetime

    $dede.Cells.Item($r, 5).Value = "mt"              # Priority
    $dede.Cells.Item($r, 8).Value = "2016-08-13 06:27:06"   # Latest Handoff Datetime

    $overview.Cells.Item($r, 7).Value = "2016-08-13 06:27:06"   # Latest HO Xliff Generate Date
}
